$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.677835822105408
$ws.Range("B1").Value = 2.968576192855835
$ws.Range("D1").Value = 0.2143935710191727
$ws.Range("E1").Value = 0.6218926310539246
